$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width ---
$ws.Columns.Item(1).ColumnWidth = 32.21875

# --- Row 1: header cell with rich text "SILAHKAN MASUKKAN`nNAMA JURUSAN" ---
$headerText = "SILAHKAN MASUKKAN`nNAMA JURUSAN"
$ws.Range("A1").Value = $headerText
$boldPart = "NAMA JURUSAN"
$startPos = $headerText.Length - $boldPart.Length + 1
$ws.Range("A1").Characters($startPos, $boldPart.Length).Font.Bold = $true

# Cell-level formatting for A1 (fill already theme 9 from template; no longer bold at cell level)
$ws.Range("A1").Font.Bold = $false
$ws.Range("A1").Font.ThemeColor = 1
$ws.Range("A1").WrapText = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Rows.Item(1).RowHeight = 50.4

# --- Rows 2 & 3: jurusan values ---
$ws.Range("A2").Value = "nama jurusan 1"
$ws.Range("A3").Value = "nama jurusan 2"

# --- Selection ---
$ws.Range("C13").Select()
